$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "started finalising the documents" - drop the "/ Admin" / ", Admin" suffix from the
# Wahlleiter/Moderator login descriptions (row 3).
# F3 is written first so the newly-appended shared-string order matches the target
# (Moderator/Wahlleiter variant before Wahlleiter/Moderator variant).
$ws.Range("F3").Value = "Login Moderator, Wahlleiter"
$ws.Range("B3").Value = "Login Wahlleiter/ Moderator"
$ws.Range("D3").Value = "Login Wahlleiter/ Moderator"

# Move the active selection to C3.
$ws.Range("C3").Select()
